# Update the "Förändrad" (Changed) date column for rows 2-16
# from serial date 45170 (2023-09-01) to 45174 (2023-09-05)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$ws.Range("C2:C16").Value = 45174
